$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Shift all timestamps forward by 2 days and replace the production values
# with the newer January production data (per commit message).
$ws.Cells.Item(2, 1).Value = 46073.01041666666
$ws.Cells.Item(2, 2).Value = 1599.067
$ws.Cells.Item(3, 1).Value = 46073.02083333334
$ws.Cells.Item(3, 2).Value = 1583.6
$ws.Cells.Item(4, 1).Value = 46073.03125
$ws.Cells.Item(4, 2).Value = 1567.819
$ws.Cells.Item(5, 1).Value = 46073.04166666666
$ws.Cells.Item(5, 2).Value = 1552.087
$ws.Cells.Item(6, 1).Value = 46073.05208333334
$ws.Cells.Item(6, 2).Value = 1527.694
$ws.Cells.Item(7, 1).Value = 46073.0625
$ws.Cells.Item(7, 2).Value = 1501.378
$ws.Cells.Item(8, 1).Value = 46073.07291666666
$ws.Cells.Item(8, 2).Value = 1476.185
$ws.Cells.Item(9, 1).Value = 46073.08333333334
$ws.Cells.Item(9, 2).Value = 1451.41
$ws.Cells.Item(10, 1).Value = 46073.09375
$ws.Cells.Item(10, 2).Value = 1402.109
$ws.Cells.Item(11, 1).Value = 46073.10416666666
$ws.Cells.Item(11, 2).Value = 1368.829
$ws.Cells.Item(12, 1).Value = 46073.11458333334
$ws.Cells.Item(12, 2).Value = 1335.304
$ws.Cells.Item(13, 1).Value = 46073.125
$ws.Cells.Item(13, 2).Value = 1301.601
$ws.Cells.Item(14, 1).Value = 46073.13541666666
$ws.Cells.Item(14, 2).Value = 1248.825
$ws.Cells.Item(15, 1).Value = 46073.14583333334
$ws.Cells.Item(15, 2).Value = 1220.87
$ws.Cells.Item(16, 1).Value = 46073.15625
$ws.Cells.Item(16, 2).Value = 1192.142
$ws.Cells.Item(17, 1).Value = 46073.16666666666
$ws.Cells.Item(17, 2).Value = 1162.905
$ws.Cells.Item(18, 1).Value = 46073.17708333334
$ws.Cells.Item(18, 2).Value = 1115.353
$ws.Cells.Item(19, 1).Value = 46073.1875
$ws.Cells.Item(19, 2).Value = 1092.907
$ws.Cells.Item(20, 1).Value = 46073.19791666666
$ws.Cells.Item(20, 2).Value = 1070.461
$ws.Cells.Item(21, 1).Value = 46073.20833333334
$ws.Cells.Item(21, 2).Value = 1047.27
$ws.Cells.Item(22, 1).Value = 46073.21875
$ws.Cells.Item(22, 2).Value = 1010.717
$ws.Cells.Item(23, 1).Value = 46073.22916666666
$ws.Cells.Item(23, 2).Value = 984.1420000000001
$ws.Cells.Item(24, 1).Value = 46073.23958333334
$ws.Cells.Item(24, 2).Value = 966.395
$ws.Cells.Item(25, 1).Value = 46073.25
$ws.Cells.Item(25, 2).Value = 939.346
$ws.Cells.Item(26, 1).Value = 46073.26041666666
$ws.Cells.Item(26, 2).Value = 901.658
$ws.Cells.Item(27, 1).Value = 46073.27083333334
$ws.Cells.Item(27, 2).Value = 875.3049999999999
$ws.Cells.Item(28, 1).Value = 46073.28125
$ws.Cells.Item(28, 2).Value = 847.294
$ws.Cells.Item(29, 1).Value = 46073.29166666666
$ws.Cells.Item(29, 2).Value = 820.127
$ws.Cells.Item(30, 1).Value = 46073.30208333334
$ws.Cells.Item(30, 2).Value = 782.283
$ws.Cells.Item(31, 1).Value = 46073.3125
$ws.Cells.Item(31, 2).Value = 756.772
$ws.Cells.Item(32, 1).Value = 46073.32291666666
$ws.Cells.Item(32, 2).Value = 734.167
$ws.Cells.Item(33, 1).Value = 46073.33333333334
$ws.Cells.Item(33, 2).Value = 712.423
$ws.Cells.Item(34, 1).Value = 46073.34375
$ws.Cells.Item(34, 2).Value = 667.518
$ws.Cells.Item(35, 1).Value = 46073.35416666666
$ws.Cells.Item(35, 2).Value = 643.956
$ws.Cells.Item(36, 1).Value = 46073.36458333334
$ws.Cells.Item(36, 2).Value = 623.4450000000001
$ws.Cells.Item(37, 1).Value = 46073.375
$ws.Cells.Item(37, 2).Value = 602.985
$ws.Cells.Item(38, 1).Value = 46073.38541666666
$ws.Cells.Item(38, 2).Value = 573.835
$ws.Cells.Item(39, 1).Value = 46073.39583333334
$ws.Cells.Item(39, 2).Value = 564.038
$ws.Cells.Item(40, 1).Value = 46073.40625
$ws.Cells.Item(40, 2).Value = 555.514
$ws.Cells.Item(41, 1).Value = 46073.41666666666
$ws.Cells.Item(41, 2).Value = 546.5069999999999
$ws.Cells.Item(42, 1).Value = 46073.42708333334
$ws.Cells.Item(42, 2).Value = 537.799
$ws.Cells.Item(43, 1).Value = 46073.4375
$ws.Cells.Item(43, 2).Value = 545.877
$ws.Cells.Item(44, 1).Value = 46073.44791666666
$ws.Cells.Item(44, 2).Value = 555.498
$ws.Cells.Item(45, 1).Value = 46073.45833333334
$ws.Cells.Item(45, 2).Value = 564.889
$ws.Cells.Item(46, 1).Value = 46073.46875
$ws.Cells.Item(46, 2).Value = 582.177
$ws.Cells.Item(47, 1).Value = 46073.47916666666
$ws.Cells.Item(47, 2).Value = 611.064
$ws.Cells.Item(48, 1).Value = 46073.48958333334
$ws.Cells.Item(48, 2).Value = 641.582
$ws.Cells.Item(49, 1).Value = 46073.5
$ws.Cells.Item(49, 2).Value = 672.793
$ws.Cells.Item(50, 1).Value = 46073.51041666666
$ws.Cells.Item(50, 2).Value = 740.553
$ws.Cells.Item(51, 1).Value = 46073.52083333334
$ws.Cells.Item(51, 2).Value = 792.776
$ws.Cells.Item(52, 1).Value = 46073.53125
$ws.Cells.Item(52, 2).Value = 846.354
$ws.Cells.Item(53, 1).Value = 46073.54166666666
$ws.Cells.Item(53, 2).Value = 900.99
$ws.Cells.Item(54, 1).Value = 46073.55208333334
$ws.Cells.Item(54, 2).Value = 1019.103
$ws.Cells.Item(55, 1).Value = 46073.5625
$ws.Cells.Item(55, 2).Value = 1110.436
$ws.Cells.Item(56, 1).Value = 46073.57291666666
$ws.Cells.Item(56, 2).Value = 1179.974
$ws.Cells.Item(57, 1).Value = 46073.58333333334
$ws.Cells.Item(57, 2).Value = 1294.811
$ws.Cells.Item(58, 1).Value = 46073.59375
$ws.Cells.Item(58, 2).Value = 1445.912
$ws.Cells.Item(59, 1).Value = 46073.60416666666
$ws.Cells.Item(59, 2).Value = 1533.49
$ws.Cells.Item(60, 1).Value = 46073.61458333334
$ws.Cells.Item(60, 2).Value = 1574.243
$ws.Cells.Item(61, 1).Value = 46073.625
$ws.Cells.Item(61, 2).Value = 1653.016
$ws.Cells.Item(62, 1).Value = 46073.63541666666
$ws.Cells.Item(62, 2).Value = 1739.951
$ws.Cells.Item(63, 1).Value = 46073.64583333334
$ws.Cells.Item(63, 2).Value = 1798.899
$ws.Cells.Item(64, 1).Value = 46073.65625
$ws.Cells.Item(64, 2).Value = 1908.164
$ws.Cells.Item(65, 1).Value = 46073.66666666666
$ws.Cells.Item(65, 2).Value = 1967.702
$ws.Cells.Item(66, 1).Value = 46073.67708333334
$ws.Cells.Item(66, 2).Value = 2057.548
$ws.Cells.Item(67, 1).Value = 46073.6875
$ws.Cells.Item(67, 2).Value = 2098.204
$ws.Cells.Item(68, 1).Value = 46073.69791666666
$ws.Cells.Item(68, 2).Value = 2138.75
$ws.Cells.Item(69, 1).Value = 46073.70833333334
$ws.Cells.Item(69, 2).Value = 2178.947
$ws.Cells.Item(70, 1).Value = 46073.71875
$ws.Cells.Item(70, 2).Value = 2229.217
$ws.Cells.Item(71, 1).Value = 46073.72916666666
$ws.Cells.Item(71, 2).Value = 2250.77
$ws.Cells.Item(72, 1).Value = 46073.73958333334
$ws.Cells.Item(72, 2).Value = 2272.106
$ws.Cells.Item(73, 1).Value = 46073.75
$ws.Cells.Item(73, 2).Value = 2292.36
$ws.Cells.Item(74, 1).Value = 46073.76041666666
$ws.Cells.Item(74, 2).Value = 2322.006
$ws.Cells.Item(75, 1).Value = 46073.77083333334
$ws.Cells.Item(75, 2).Value = 2328.441
$ws.Cells.Item(76, 1).Value = 46073.78125
$ws.Cells.Item(76, 2).Value = 2334.06
$ws.Cells.Item(77, 1).Value = 46073.79166666666
$ws.Cells.Item(77, 2).Value = 2338.934
$ws.Cells.Item(78, 1).Value = 46073.80208333334
$ws.Cells.Item(78, 2).Value = 2340.843
$ws.Cells.Item(79, 1).Value = 46073.8125
$ws.Cells.Item(79, 2).Value = 2339.701
$ws.Cells.Item(80, 1).Value = 46073.82291666666
$ws.Cells.Item(80, 2).Value = 2339.241
$ws.Cells.Item(81, 1).Value = 46073.83333333334
$ws.Cells.Item(81, 2).Value = 2338.972
$ws.Cells.Item(82, 1).Value = 46073.84375
$ws.Cells.Item(82, 2).Value = 2343.798
$ws.Cells.Item(83, 1).Value = 46073.85416666666
$ws.Cells.Item(83, 2).Value = 2334.354
$ws.Cells.Item(84, 1).Value = 46073.86458333334
$ws.Cells.Item(84, 2).Value = 2324.34
$ws.Cells.Item(85, 1).Value = 46073.875
$ws.Cells.Item(85, 2).Value = 2252.323
$ws.Cells.Item(86, 1).Value = 46073.88541666666
$ws.Cells.Item(86, 2).Value = 2305.569
$ws.Cells.Item(87, 1).Value = 46073.89583333334
$ws.Cells.Item(87, 2).Value = 2242.223
$ws.Cells.Item(88, 1).Value = 46073.90625
$ws.Cells.Item(88, 2).Value = 2238.162
$ws.Cells.Item(89, 1).Value = 46073.91666666666
$ws.Cells.Item(89, 2).Value = 2236.393
$ws.Cells.Item(90, 1).Value = 46073.92708333334
$ws.Cells.Item(90, 2).Value = 2225.751
$ws.Cells.Item(91, 1).Value = 46073.9375
$ws.Cells.Item(91, 2).Value = 2220.989
$ws.Cells.Item(92, 1).Value = 46073.94791666666
$ws.Cells.Item(92, 2).Value = 2216.946
$ws.Cells.Item(93, 1).Value = 46073.95833333334
$ws.Cells.Item(93, 2).Value = 2212.739
$ws.Cells.Item(94, 1).Value = 46073.96875
$ws.Cells.Item(94, 2).Value = 0
$ws.Cells.Item(95, 1).Value = 46073.97916666666
$ws.Cells.Item(95, 2).Value = 0
$ws.Cells.Item(96, 1).Value = 46073.98958333334
$ws.Cells.Item(96, 2).Value = 0
$ws.Cells.Item(97, 1).Value = 46074
$ws.Cells.Item(97, 2).Value = 0
